$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures to several
# rows across multiple sheets (ALC, BSM, CRP, GSM, LTW, WVR), as produced
# by the scheduled market-data refresh. Columns H-N hold:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 570.70734
$ws.Range("J17").Value = 577.55
$ws.Range("L17").Value = 1732.65
$ws.Range("N17").Value = -2068.65

$ws.Range("H86").Value = 7232.9165
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 9143.888999999999
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 9143.888999999999
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -11389.889

$ws.Range("H89").Value = 7232.9165
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 9143.888999999999
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 45719.44499999999
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -56951.44499999999

$ws.Range("H113").Value = 64293.75
$ws.Range("I113").Value = 168218.33
$ws.Range("J113").Value = 1939
$ws.Range("K113").Value = 168218.33
$ws.Range("L113").Value = 1939
$ws.Range("M113").Value = -164964.33
$ws.Range("N113").Value = -8447

$ws.Range("H128").Value = 45379.668
$ws.Range("J128").Value = 45379.668
$ws.Range("L128").Value = 45379.668
$ws.Range("N128").Value = -55339.668

$ws.Range("H129").Value = 941.625
$ws.Range("I129").Value = 534.75
$ws.Range("J129").Value = 1023
$ws.Range("K129").Value = 1604.25
$ws.Range("L129").Value = 3069
$ws.Range("M129").Value = 3395.75
$ws.Range("N129").Value = -13069

$ws.Range("H138").Value = 7212.4043
$ws.Range("I138").Value = 1489.3704
$ws.Range("J138").Value = 14938.5
$ws.Range("K138").Value = 4468.1112
$ws.Range("L138").Value = 44815.5
$ws.Range("M138").Value = 671.8887999999997
$ws.Range("N138").Value = -55095.5

$ws.Range("H141").Value = 2960.238
$ws.Range("I141").Value = 2715
$ws.Range("J141").Value = 4002.5
$ws.Range("K141").Value = 8145
$ws.Range("L141").Value = 12007.5
$ws.Range("M141").Value = -2965
$ws.Range("N141").Value = -22367.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346

$ws.Range("H69").Value = 42995
$ws.Range("J69").Value = 42995
$ws.Range("L69").Value = 42995
$ws.Range("N69").Value = -44617

$ws.Range("H72").Value = 42995
$ws.Range("J72").Value = 42995
$ws.Range("L72").Value = 128985
$ws.Range("N72").Value = -137097

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37757.227
$ws.Range("I31").Value = 1100.8
$ws.Range("K31").Value = 1100.8
$ws.Range("M31").Value = -805.8

$ws.Range("H34").Value = 37757.227
$ws.Range("I34").Value = 1100.8
$ws.Range("K34").Value = 1100.8
$ws.Range("M34").Value = -898.8

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5863571.5
$ws.Range("I11").Value = 6237692.5
$ws.Range("K11").Value = 6237692.5
$ws.Range("M11").Value = -6237553.5

$ws.Range("H18").Value = 2500
$ws.Range("J18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("N18").Value = -3086

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

$ws.Range("H62").Value = 11998
$ws.Range("J62").Value = 11998
$ws.Range("L62").Value = 11998
$ws.Range("N62").Value = -13370

$ws.Range("H65").Value = 11998
$ws.Range("J65").Value = 11998
$ws.Range("L65").Value = 35994
$ws.Range("N65").Value = -42858

$ws.Range("H70").Value = 158378.39
$ws.Range("I70").Value = 289315.16
$ws.Range("J70").Value = 5618.8335
$ws.Range("K70").Value = 289315.16
$ws.Range("L70").Value = 5618.8335
$ws.Range("M70").Value = -289045.16
$ws.Range("N70").Value = -6158.8335

$ws.Range("H73").Value = 158378.39
$ws.Range("I73").Value = 289315.16
$ws.Range("J73").Value = 5618.8335
$ws.Range("K73").Value = 289315.16
$ws.Range("L73").Value = 5618.8335
$ws.Range("M73").Value = -288379.16
$ws.Range("N73").Value = -7490.8335

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5183.1665
$ws.Range("I7").Value = 1999.6666
$ws.Range("J7").Value = 8366.666999999999
$ws.Range("K7").Value = 1999.6666
$ws.Range("L7").Value = 8366.666999999999
$ws.Range("M7").Value = -1887.6666
$ws.Range("N7").Value = -8590.666999999999

$ws.Range("H22").Value = 1384.08
$ws.Range("I22").Value = 2649.6667
$ws.Range("J22").Value = 984.4211
$ws.Range("K22").Value = 2649.6667
$ws.Range("L22").Value = 984.4211
$ws.Range("M22").Value = -2354.6667
$ws.Range("N22").Value = -1574.4211

$ws.Range("H23").Value = 40003.5
$ws.Range("J23").Value = 40003.5
$ws.Range("L23").Value = 40003.5
$ws.Range("N23").Value = -40463.5

$ws.Range("H27").Value = 1384.08
$ws.Range("I27").Value = 2649.6667
$ws.Range("J27").Value = 984.4211
$ws.Range("K27").Value = 2649.6667
$ws.Range("L27").Value = 984.4211
$ws.Range("M27").Value = -2542.6667
$ws.Range("N27").Value = -1198.4211

$ws.Range("H68").Value = 3026.4666
$ws.Range("I68").Value = 1516.6666
$ws.Range("J68").Value = 4033
$ws.Range("K68").Value = 1516.6666
$ws.Range("L68").Value = 4033
$ws.Range("M68").Value = -767.6666
$ws.Range("N68").Value = -5531

$ws.Range("H71").Value = 3026.4666
$ws.Range("I71").Value = 1516.6666
$ws.Range("J71").Value = 4033
$ws.Range("K71").Value = 7583.333000000001
$ws.Range("L71").Value = 20165
$ws.Range("M71").Value = -3839.333000000001
$ws.Range("N71").Value = -27653

$ws.Range("H126").Value = 5183.1665
$ws.Range("I126").Value = 1999.6666
$ws.Range("J126").Value = 8366.666999999999
$ws.Range("K126").Value = 5998.9998
$ws.Range("L126").Value = 25100.001
$ws.Range("M126").Value = -3528.9998
$ws.Range("N126").Value = -30040.001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

$ws.Range("H132").Value = 1573.0507
$ws.Range("I132").Value = 1452.1343
$ws.Range("J132").Value = 2248.1667
$ws.Range("K132").Value = 4356.4029
$ws.Range("L132").Value = 6744.500100000001
$ws.Range("M132").Value = -1826.4029
$ws.Range("N132").Value = -11804.5001

